# Removed semicolon from initial .xlsx files
# Strips the trailing ';' from each email address displayed in column B
# (rows 2-12). The underlying mailto: hyperlink targets are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $email = $cell.Value()
    if ($email -ne $null -and $email.EndsWith(";")) {
        $cell.Value = $email.TrimEnd(";")
    }
}
